$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new test case "TC 2" with data "Dhaka"
$ws.Range("A3").Value = "TC 2"
$ws.Range("B3").Value = "Dhaka"

# Row 4: new test case "TC 3" with data "USA"
# (value entry order matches shared-string insertion order: B4 before A4)
$ws.Range("B4").Value = "USA"
$ws.Range("A4").Value = "TC 3"

# Update the current selection to match the final cursor position
$ws.Range("A4").Select()
